$d = $word.ActiveDocument

# Locate the "Data Engineering and Infrastructure Architecture" paragraph under the
# Siege Analytics / PARTNER entry (it is the unique occurrence in the document).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Data Engineering and Infrastructure Architecture") {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    Write-Output "Target paragraph not found"
} else {
    $newLines = @(
        "• Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections",
        "• Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government",
        "• Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations"
    )

    $current = $targetPara
    foreach ($line in $newLines) {
        $current.Range.InsertParagraphAfter()
        $current = $current.Next()
        $current.Range.Text = $line
    }

    Write-Output "Inserted $($newLines.Count) paragraphs after 'Data Engineering and Infrastructure Architecture'"
}
